# feat: add 2022-Q3 data
#
# Starting layout:  Sheet1 = "总计" (summary), Sheet2 = "2021-Q3" (fund
# holdings for that quarter).
#
# Target layout:    Sheet1 = "总计", Sheet2 = "2022-Q3" (new fund
# holdings), Sheet3 = "2021-Q3" (old fund holdings, unchanged, just moved
# to a new sheet).

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Summary sheet ("总计"): insert a new top data row with the 2022-Q3
#    totals and push the existing 2021-Q3 totals down to row 3.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Carry row 2's formatting down to row 3 before overwriting row 2, so the
# "序号" cell keeps its centered/bold/bordered look.
$summary.Range("A2").Copy($summary.Range("A3"))
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.08

# New 2022-Q3 summary row (A2 already has the value/format we need: 0).
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.11

# ---------------------------------------------------------------------
# 2. Detail sheets. The existing "2021-Q3" holdings sheet is kept (so its
#    underlying worksheet part / tab formatting survives) but renamed to
#    "2022-Q3" and refilled with the new quarter's holdings. A brand-new
#    sheet named "2021-Q3" is added to hold the data that used to live on
#    the original sheet.
# ---------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item(2)

# Add the replacement "2021-Q3" sheet right after the current one.
$movedQ3 = $wb.Worksheets.Add($null, $oldQ3)

# Clone the header-row / row-label formatting (style only) from the
# still-untouched original sheet before its content changes.
$oldQ3.Range("B1:H1").Copy()
$movedQ3.Range("B1:H1").PasteSpecial($xlPasteFormats)
$oldQ3.Range("A2").Copy()
$movedQ3.Range("A2").PasteSpecial($xlPasteFormats)

# Recreate the (unchanged) 2021-Q3 holdings table on the new sheet.
$movedQ3.Range("B1").Value = "基金代码"
$movedQ3.Range("C1").Value = "基金名称"
$movedQ3.Range("D1").Value = "基金金额"
$movedQ3.Range("E1").Value = "股票总仓位"
$movedQ3.Range("F1").Value = "仓位占比"
$movedQ3.Range("G1").Value = "持有市值(亿元)"
$movedQ3.Range("H1").Value = "仓位排名"

$movedQ3.Range("A2").Value = 0
$movedQ3.Range("B2").Value = "'001637"
$movedQ3.Range("B2").Style = "Normal"
$movedQ3.Range("C2").Value = "嘉实腾讯自选股大数据策略股票"
$movedQ3.Range("D2").Value = "'4.81"
$movedQ3.Range("D2").Style = "Normal"
$movedQ3.Range("E2").Value = "'88.55"
$movedQ3.Range("E2").Style = "Normal"
$movedQ3.Range("F2").Value = "'1.75"
$movedQ3.Range("F2").Style = "Normal"
$movedQ3.Range("G2").Value = "'0.0842"
$movedQ3.Range("G2").Style = "Normal"
$movedQ3.Range("H2").Value = 5

# Rename: the original sheet becomes "2022-Q3"; the new sheet takes over
# the "2021-Q3" name.
$oldQ3.Name = "2022-Q3"
$movedQ3.Name = "2021-Q3"

# ---------------------------------------------------------------------
# 3. Populate "2022-Q3" (the renamed original sheet) with the new
#    quarter's fund holdings.
# ---------------------------------------------------------------------
$newQ3 = $oldQ3

$newQ3.Range("B1").Value = "基金代码"
$newQ3.Range("C1").Value = "基金名称"
$newQ3.Range("D1").Value = "基金规模"
$newQ3.Range("E1").Value = "股票总仓位"
$newQ3.Range("F1").Value = "仓位占比"
$newQ3.Range("G1").Value = "持有市值(亿元)"
$newQ3.Range("H1").Value = "仓位排名"

$newQ3.Range("A2").Value = 0
$newQ3.Range("B2").Value = "'015784"
$newQ3.Range("B2").Style = "Normal"
$newQ3.Range("C2").Value = "中信建投中证1000指数增强A"
$newQ3.Range("D2").Value = "'8.10"
$newQ3.Range("D2").Style = "Normal"
$newQ3.Range("E2").Value = "'92.20"
$newQ3.Range("E2").Style = "Normal"
$newQ3.Range("F2").Value = "'0.69"
$newQ3.Range("F2").Style = "Normal"
$newQ3.Range("G2").Value = "'0.0559"
$newQ3.Range("G2").Style = "Normal"
$newQ3.Range("H2").Value = 3

$newQ3.Range("A3").Value = 1
$newQ3.Range("B3").Value = "'013466"
$newQ3.Range("B3").Style = "Normal"
$newQ3.Range("C3").Value = "博时智选量化多因子股票C"
$newQ3.Range("D3").Value = "'2.28"
$newQ3.Range("D3").Style = "Normal"
$newQ3.Range("E3").Value = "'92.38"
$newQ3.Range("E3").Style = "Normal"
$newQ3.Range("F3").Value = "'1.09"
$newQ3.Range("F3").Style = "Normal"
$newQ3.Range("G3").Value = "'0.0249"
$newQ3.Range("G3").Style = "Normal"
$newQ3.Range("H3").Value = 7

$newQ3.Range("A4").Value = 2
$newQ3.Range("B4").Value = "'015785"
$newQ3.Range("B4").Style = "Normal"
$newQ3.Range("C4").Value = "中信建投中证1000指数增强C"
$newQ3.Range("D4").Value = "'3.32"
$newQ3.Range("D4").Style = "Normal"
$newQ3.Range("E4").Value = "'92.20"
$newQ3.Range("E4").Style = "Normal"
$newQ3.Range("F4").Value = "'0.69"
$newQ3.Range("F4").Style = "Normal"
$newQ3.Range("G4").Value = "'0.0229"
$newQ3.Range("G4").Style = "Normal"
$newQ3.Range("H4").Value = 3

$newQ3.Range("A5").Value = 3
$newQ3.Range("B5").Value = "'013465"
$newQ3.Range("B5").Style = "Normal"
$newQ3.Range("C5").Value = "博时智选量化多因子股票A"
$newQ3.Range("D5").Value = "'0.49"
$newQ3.Range("D5").Style = "Normal"
$newQ3.Range("E5").Value = "'92.38"
$newQ3.Range("E5").Style = "Normal"
$newQ3.Range("F5").Value = "'1.09"
$newQ3.Range("F5").Style = "Normal"
$newQ3.Range("G5").Value = "'0.0053"
$newQ3.Range("G5").Style = "Normal"
$newQ3.Range("H5").Value = 7

# Apply the "总计" header/row-label style (bold, centered, bordered) to
# the new sheet's header row and A-column, matching the rest of the
# workbook (same style index used by the "2022-Q3" sheet in the target).
$summary.Range("B1").Copy()
$newQ3.Range("B1:H1").PasteSpecial($xlPasteFormats)
$summary.Range("A2").Copy()
$newQ3.Range("A2:A5").PasteSpecial($xlPasteFormats)

# Leave the workbook the way it started: "总计" as the active sheet.
$summary.Activate()
